$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new values (B11 "Startet kl. 12:00" and B13 "74 hours")
$ws.Range("B11").Value = "Startet kl. 12:00"
$ws.Range("B13").Value = "74 hours"

# Update the active cell selection to B14
$ws.Range("B14").Select()
